# "multiple_testing_correction stond 2x in de sheet"
# The "algemeen" sheet had the setting "multiple_testing_correction" twice
# (column I, value "BH", and column S, value "bonferroni"). Remove the
# duplicate (first) occurrence in column I so the setting exists only once.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("algemeen")

# Column I currently holds the duplicate header "multiple_testing_correction"
# with value "BH". Deleting the whole column shifts every later column
# (J..T) one position to the left.
$ws.Columns.Item(9).Delete()

# Re-select "algemeen" as the active sheet/cell, matching the saved view
# state after the edit.
$ws.Activate()
$ws.Range("L8").Select()
